# Update the division problems in the table to the newly generated values.
$d = $word.ActiveDocument

$replacements = @(
    @("746÷8=", "374÷2="),
    @("985÷3=", "161÷3="),
    @("893÷8=", "213÷6="),
    @("445÷9=", "846÷9="),
    @("612÷8=", "297÷5="),
    @("272÷9=", "337÷5="),
    @("588÷2=", "664÷4="),
    @("385÷5=", "756÷9="),
    @("276÷6=", "637÷4="),
    @("303÷2=", "407÷2="),
    @("550÷3=", "722÷3="),
    @("353÷3=", "390÷7="),
    @("137÷2=", "227÷6="),
    @("851÷7=", "198÷6="),
    @("649÷8=", "242÷9="),
    @("183÷9=", "266÷9="),
    @("916÷6=", "134÷9="),
    @("386÷8=", "194÷9="),
    @("110÷9=", "793÷7="),
    @("126÷5=", "730÷2="),
    @("688÷2=", "255÷2="),
    @("750÷7=", "826÷9="),
    @("876÷6=", "897÷9="),
    @("508÷9=", "529÷3="),
    @("216÷6=", "173÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
